$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new rows to the Parameter / Value / Description table ---
# Row 56: MFTC_WEP_scaling = 1
# Row 57: WFF_or_Benefit   = Max
#
# The "'" prefix forces Excel to store the numeric-looking "1" as text,
# matching how every other Value-column entry in this sheet is stored
# (all of them are shared-string / text cells, even numeric-looking ones).

$ws.Range("B56").Value2 = "MFTC_WEP_scaling"
$ws.Range("C56").Value2 = "'1"
$ws.Range("D56").Value2 = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

$ws.Range("B57").Value2 = "WFF_or_Benefit"
$ws.Range("C57").Value2 = "Max"
$ws.Range("D57").Value2 = "What work decision should we assume? Go off-benefit and receive IWTC = ""WFF"", stay on-benefit = ""Benefit"", or whichever gives a higher net income = ""Max"""

# --- Formatting: mirror the look of the last existing data row (55) ---
# Copy font / fill / alignment from row 55, then strip the top border that
# row 55 has (row 55 is the bottom border of the previous group; the new
# rows are a fresh, borderless group of their own).
$src = $ws.Range("B55:D55")
$dst = $ws.Range("B56:D57")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$dst.Borders.LineStyle = -4142

Write-Output "Added MFTC_WEP_scaling (row 56) and WFF_or_Benefit (row 57) parameters."
